$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C16").Value = "Test dataset - Duplicate value in FA19 column.xlsx"
$ws.Range("D16").Value = "\Testdata\Templates\ImportPublications\Staging_Env\Test dataset - Duplicate value in FA19 column.xlsx"

$ws.Range("D16").Select()
